$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Change 1: heading text "Divieto di prosecuzione dell'attivita" ->
#           "Misure urgenti di messa in sicurezza"
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "Divieto di prosecuzione dell’attività", $true, $false, $false, $false, $false,
    $true, 1, $false, "Misure urgenti di messa in sicurezza", 2) | Out-Null

# --------------------------------------------------------------------
# Change 2: split the single "Pertanto questo Comando diffida..."
# paragraph into four paragraphs:
#   [Senza 758]
#   Pertanto questo Comando **diffida** il responsabile ... 151/2011.
#   [Con 758]
#   Inoltre questo Comando comunica ... **eliminare** le carenza ... .
# --------------------------------------------------------------------

# Locate the target paragraph by its text.
$targetIndex = -1
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Pertanto questo Comando diffida*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs($targetIndex)
$targetRange = $target.Range
# Keep the original wording (with its original typography, e.g. the
# non-breaking spaces already present before "4" and "151/2011") so the
# reconstructed runs reuse exactly the same characters.
$originalFull = $d.Range($targetRange.Start, $targetRange.End - 1)
$originalText = $originalFull.Text

# Insert a blank paragraph before the target ("[Senza 758]" slot), and two
# blank paragraphs after it ("[Con 758]" and the "Inoltre ..." slots), so
# we end up with four paragraph slots in document order.
$ins = $targetRange.Duplicate
$ins.Collapse(1)
$ins.InsertParagraphBefore()

$targetIndex = $targetIndex + 1
$target = $d.Paragraphs($targetIndex)

$insAfter = $target.Range.Duplicate
$insAfter.Collapse(0)
$insAfter.InsertParagraphAfter()
$insAfter2 = $d.Paragraphs($targetIndex + 1).Range.Duplicate
$insAfter2.Collapse(0)
$insAfter2.InsertParagraphAfter()

# Paragraph indices after insertion:
$senzaIndex   = $targetIndex - 1
$mainIndex    = $targetIndex
$conIndex     = $targetIndex + 1
$inoltreIndex = $targetIndex + 2

function Set-ParaText($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $r.Text = $text
}

function Style-ParaRun($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $r.Style = "Da modificare EG"
}

function Bold-Substring($paraIndex, $needle) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $text = $d.Range($full.Start, $full.End - 1).Text
    $relStart = $text.IndexOf($needle)
    $start = $full.Start + $relStart
    $end = $start + $needle.Length
    $r = $d.Range($start, $end)
    $r.Bold = 1
}

# --- [Senza 758] ---
Set-ParaText $senzaIndex "[Senza 758]"
Style-ParaRun $senzaIndex

# --- main paragraph (reuse the original wording verbatim) ---
Set-ParaText $mainIndex $originalText
Style-ParaRun $mainIndex
Bold-Substring $mainIndex "diffida"

# --- [Con 758] ---
Set-ParaText $conIndex "[Con 758]"
Style-ParaRun $conIndex

# --- Inoltre questo Comando ... ---
$nbsp = [char]0x00A0
$inoltreText = "Inoltre questo Comando comunica che ha attivato procedura sanzionatoria ai sensi dell’articolo" + $nbsp + "20 del decreto legislativo" + $nbsp + "758/94, al fine di eliminare le carenza sopra indicate."
Set-ParaText $inoltreIndex $inoltreText
Style-ParaRun $inoltreIndex
Bold-Substring $inoltreIndex "eliminare"

# --------------------------------------------------------------------
# Change 3: Normal style overflowPunct true -> false
# (exposed on the Word OM as ParagraphFormat.HangingPunctuation)
# --------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $false
